$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential-notice date from 2021-05-12 to 2021-05-13
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."
$ws.Rows.Item(16).AutoFit()

# Updated Weight (D) and Percent Change (E) values for rows 2-13
$ws.Cells.Item(2, 4).Value2 = 0.03077848746848751
$ws.Cells.Item(2, 5).Value2 = 0.01411930815390061
$ws.Cells.Item(3, 4).Value2 = 0.02384427320938459
$ws.Cells.Item(3, 5).Value2 = 0.02045235803657341
$ws.Cells.Item(4, 4).Value2 = 0.05280549727367765
$ws.Cells.Item(4, 5).Value2 = 0.01461919358641839
$ws.Cells.Item(5, 4).Value2 = 0.1392815760960588
$ws.Cells.Item(5, 5).Value2 = 0.007619678648335126
$ws.Cells.Item(6, 4).Value2 = 0.03243742544544743
$ws.Cells.Item(6, 5).Value2 = -0.01407459535538358
$ws.Cells.Item(7, 4).Value2 = 0.1171713142157923
$ws.Cells.Item(7, 5).Value2 = 0.008864779042074611
$ws.Cells.Item(8, 4).Value2 = 0.1028968792434755
$ws.Cells.Item(8, 5).Value2 = 0.01893158388003724
$ws.Cells.Item(9, 4).Value2 = 0.02982682591820203
$ws.Cells.Item(9, 5).Value2 = 0.01607180129409302
$ws.Cells.Item(10, 4).Value2 = 0.1271974519493925
$ws.Cells.Item(10, 5).Value2 = 0.02092130518234159
$ws.Cells.Item(11, 4).Value2 = 0.241118517263453
$ws.Cells.Item(11, 5).Value2 = 0.01067718353212777
$ws.Cells.Item(12, 4).Value2 = 0.1026417519166288
$ws.Cells.Item(12, 5).Value2 = 0.00707292707292706
$ws.Cells.Item(13, 4).Value2 = 1
$ws.Cells.Item(13, 5).Value2 = 0.01172660409170523

$ws.Protect()
